$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (old Amana Takaful row; its data moves into row 6 with recalculated values)
$ws.Rows("7").Delete()

# --- Row 2 ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("B2").ClearFormats()
$ws.Range("D2").Value = 0.12825
$ws.Range("E2").Value = 0.175
$ws.Range("G2").Value = 0.1910355352081971
$ws.Range("H2").Value = 0.1910355352081971
$ws.Range("I2").Value = 0.1594419010246348
$ws.Range("J2").Value = 0.117971952515715
$ws.Range("K2").Value = 58.38099999999999
$ws.Range("L2").Value = 0.1272749073468498
$ws.Range("M2").Value = 7.382
$ws.Range("N2").Value = 0.01953943885653785
$ws.Range("O2").Value = 0.1264452475976773
$ws.Range("P2").Value = 7.382
$ws.Range("Q2").Value = 0.01953943885653785
$ws.Range("R2").Value = 0.1264452475976773
$ws.Range("U2").Value = 13.789
$ws.Range("V2").Value = 0.03649814716781366
$ws.Range("W2").Value = 0.158333622604832
$ws.Range("X2").Value = 0.09141942577696022
$ws.Range("Y2").Value = 0.06691419682787175
$ws.Range("Z2").Value = 1.244762607732888
$ws.Range("AA2").Value = 0.16441082774004
$ws.Range("AB2").Value = 0.08995604703923216
$ws.Range("AC2").Value = 0.07498937632819515
$ws.Range("AD2").Value = 17.76
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 17.76
$ws.Range("AG2").Value = 3.971
$ws.Range("AH2").Value = 0.0448983719284053
$ws.Range("AI2").Value = 0.03994961310059385
$ws.Range("AJ2").Value = 0.01040152342634721
$ws.Range("AK2").Value = 0.009218354996041981
$ws.Range("AL2").Value = 2.862
$ws.Range("AM2").Value = 2.862
$ws.Range("AN2").Value = 0.224098118635727
$ws.Range("AO2").Value = 25.55415793151642
$ws.Range("AP2").Value = 0.05010662326027432
$ws.Range("AQ2").Value = 25.55415793151642

# --- Row 3 ---
$ws.Range("D3").Value = 0.157
$ws.Range("E3").Value = 0.429
$ws.Range("G3").Value = 0.1392282958199357
$ws.Range("H3").Value = 0.1392282958199357
$ws.Range("I3").Value = 0.1559485530546623
$ws.Range("J3").Value = 0.112384663777436
$ws.Range("K3").Value = 6.64
$ws.Range("L3").Value = 0.1067524115755627
$ws.Range("M3").Value = 1.9
$ws.Range("N3").Value = 0.04042553191489361
$ws.Range("O3").Value = 0.286144578313253
$ws.Range("P3").Value = 1.9
$ws.Range("Q3").Value = 0.04042553191489361
$ws.Range("R3").Value = 0.286144578313253
$ws.Range("U3").Value = 3.15
$ws.Range("V3").Value = 0.06702127659574468
$ws.Range("W3").Value = 0.1964497041420118
$ws.Range("X3").Value = 0.09150547839753662
$ws.Range("Y3").Value = 0.1049442257444752
$ws.Range("Z3").Value = 1.857825567502987
$ws.Range("AA3").Value = 0.2087911017609475
$ws.Range("AB3").Value = 0.08999707598343056
$ws.Range("AC3").Value = 0.118794025777517
$ws.Range("AD3").Value = 2.49
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 2.49
$ws.Range("AG3").Value = -0.6599999999999997
$ws.Range("AH3").Value = 0.0503131945847646
$ws.Range("AI3").Value = 0.05546892403653375
$ws.Range("AJ3").Value = -0.01424255502805351
$ws.Range("AK3").Value = -0.01581217057977958
$ws.Range("AL3").Value = 0.506
$ws.Range("AM3").Value = 0.506
$ws.Range("AN3").Value = 0.249498997995992
$ws.Range("AO3").Value = 19.1699604743083
$ws.Range("AP3").Value = -0.06613226452905809
$ws.Range("AQ3").Value = 19.1699604743083

# --- Row 4 ---
$ws.Range("T4").ClearContents()
$ws.Range("D4").Value = 0.277
$ws.Range("E4").Value = 0.146
$ws.Range("G4").Value = 0.1955445544554456
$ws.Range("H4").Value = 0.1955445544554456
$ws.Range("I4").Value = 0.09331683168316832
$ws.Range("J4").Value = 0.08733894141293828
$ws.Range("K4").Value = 6.22
$ws.Range("L4").Value = 0.07698019801980198
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 0.02307692307692308
$ws.Range("W4").Value = 0.1495192307692308
$ws.Range("X4").Value = 0.08933775359786794
$ws.Range("Y4").Value = 0.06018147717136282
$ws.Range("Z4").Value = 1.984819081775529
$ws.Range("AA4").Value = 0.1733519974984749
$ws.Range("AB4").Value = 0.0889278847286559
$ws.Range("AC4").Value = 0.08442411276981901
$ws.Range("AD4").Value = 1.05
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.05
$ws.Range("AG4").Value = -0.5700000000000001
$ws.Range("AH4").Value = 0.01473684210526316
$ws.Range("AI4").Value = 0.02145045965270685
$ws.Range("AJ4").Value = -0.008186126669538992
$ws.Range("AK4").Value = -0.01204310162687513
$ws.Range("AL4").Value = 0.889
$ws.Range("AM4").Value = 0.889
$ws.Range("AN4").Value = 0.1265060240963855
$ws.Range("AO4").Value = 8.481439820022498
$ws.Range("AP4").Value = -0.06867469879518072
$ws.Range("AQ4").Value = 8.481439820022498

# --- Row 5 ---
$ws.Range("D5").Value = 0.09949999999999999
$ws.Range("E5").Value = 0.175
$ws.Range("G5").Value = 0.2148397976391231
$ws.Range("H5").Value = 0.2148397976391231
$ws.Range("I5").Value = 0.1865092748735244
$ws.Range("J5").Value = 0.1497725995196484
$ws.Range("K5").Value = 46.3
$ws.Range("L5").Value = 0.1561551433389545
$ws.Range("M5").Value = 5.39
$ws.Range("N5").Value = 0.02154276578737011
$ws.Range("O5").Value = 0.116414686825054
$ws.Range("P5").Value = 5.39
$ws.Range("Q5").Value = 0.02154276578737011
$ws.Range("R5").Value = 0.116414686825054
$ws.Range("U5").Value = 8.380000000000001
$ws.Range("V5").Value = 0.03349320543565149
$ws.Range("W5").Value = 0.1671480144404332
$ws.Range("X5").Value = 0.0913333731563838
$ws.Range("Y5").Value = 0.0758146412840494
$ws.Range("Z5").Value = 1.038038055560418
$ws.Range("AA5").Value = 0.155469657981605
$ws.Range("AB5").Value = 0.08991501809503374
$ws.Range("AC5").Value = 0.06555463988657129
$ws.Range("AD5").Value = 12.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 12.5
$ws.Range("AG5").Value = 4.119999999999999
$ws.Range("AH5").Value = 0.0475827940616673
$ws.Range("AI5").Value = 0.03699319325244155
$ws.Range("AJ5").Value = 0.01620006291286568
$ws.Range("AK5").Value = 0.01250303471716436
$ws.Range("AL5").Value = 1.27
$ws.Range("AM5").Value = 1.27
$ws.Range("AN5").Value = 0.2072968490878939
$ws.Range("AO5").Value = 43.54330708661417
$ws.Range("AP5").Value = 0.06832504145936981
$ws.Range("AQ5").Value = 43.54330708661417

# --- Row 6 ---
$ws.Range("E6").ClearContents()
$ws.Range("B6").Value = "Amãna Takaful PLC (COSE:ATL.N0000)"
$ws.Range("D6").Value = 0.06150000000000001
$ws.Range("G6").Value = -0.02770833333333333
$ws.Range("H6").Value = -0.02770833333333333
$ws.Range("I6").Value = 0.03104166666666667
$ws.Range("J6").Value = 0.01552083333333333
$ws.Range("K6").Value = -0.779
$ws.Range("L6").Value = -0.04057291666666667
$ws.Range("M6").Value = 0.092
$ws.Range("N6").Value = 0.008846153846153846
$ws.Range("O6").Value = -0.1181001283697047
$ws.Range("P6").Value = 0.092
$ws.Range("Q6").Value = 0.008846153846153846
$ws.Range("R6").Value = -0.1181001283697047
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0.639
$ws.Range("V6").Value = 0.06144230769230769
$ws.Range("W6").Value = -0.09251781472684086
$ws.Range("X6").Value = 0.09791409736238833
$ws.Range("Y6").Value = -0.1904319120892292
$ws.Range("Z6").Value = 2.211981566820277
$ws.Range("AA6").Value = 0.03433179723502304
$ws.Range("AB6").Value = 0.09385148089526124
$ws.Range("AC6").Value = -0.0595196836602382
$ws.Range("AD6").Value = 1.72
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 1.72
$ws.Range("AG6").Value = 1.081
$ws.Range("AH6").Value = 0.1419141914191419
$ws.Range("AI6").Value = 0.1341653666146646
$ws.Range("AJ6").Value = 0.09415556136225067
$ws.Range("AK6").Value = 0.08874476643953699
$ws.Range("AL6").Value = 0.197
$ws.Range("AM6").Value = 0.197
$ws.Range("AN6").Value = 2.563338301043219
$ws.Range("AO6").Value = 3.025380710659898
$ws.Range("AP6").Value = 1.611028315946349
$ws.Range("AQ6").Value = 3.025380710659898

